$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 253.16667
$ws.Range("I135").Value = 179.29268
$ws.Range("J135").Value = 685.8570999999999
$ws.Range("K135").Value = 1613.63412
$ws.Range("L135").Value = 6172.7139
$ws.Range("M135").Value = 921.3658800000001
$ws.Range("N135").Value = -11242.7139
$ws.Range("H137").Value = 1755.9524
$ws.Range("I137").Value = 1320.898
$ws.Range("J137").Value = 3278.6428
$ws.Range("K137").Value = 3962.694
$ws.Range("L137").Value = 9835.928400000001
$ws.Range("M137").Value = -1412.694
$ws.Range("N137").Value = -14935.9284
$ws.Range("H138").Value = 2992.5293
$ws.Range("I138").Value = 1003.89655
$ws.Range("J138").Value = 4471.2563
$ws.Range("K138").Value = 3011.68965
$ws.Range("L138").Value = 13413.7689
$ws.Range("M138").Value = 2128.31035
$ws.Range("N138").Value = -23693.7689
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17330.648
$ws.Range("I32").Value = 18874.754
$ws.Range("J32").Value = 8966.75
$ws.Range("K32").Value = 18874.754
$ws.Range("L32").Value = 8966.75
$ws.Range("M32").Value = -18587.754
$ws.Range("N32").Value = -9540.75
$ws.Range("H122").Value = 3125624
$ws.Range("I122").Value = 639.9487
$ws.Range("K122").Value = 1919.8461
$ws.Range("M122").Value = 530.1538999999998
$ws.Range("H132").Value = 3825.074
$ws.Range("I132").Value = 1255.4147
$ws.Range("J132").Value = 11929.385
$ws.Range("K132").Value = 3766.2441
$ws.Range("L132").Value = 35788.155
$ws.Range("M132").Value = -1236.2441
$ws.Range("N132").Value = -40848.155
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 39126.82
$ws.Range("I134").Value = 3522.2632
$ws.Range("J134").Value = 114292
$ws.Range("K134").Value = 10566.7896
$ws.Range("L134").Value = 342876
$ws.Range("M134").Value = -8031.7896
$ws.Range("N134").Value = -347946
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2972.1184
$ws.Range("I31").Value = 2883.3276
$ws.Range("J31").Value = 3258.2222
$ws.Range("K31").Value = 2883.3276
$ws.Range("L31").Value = 3258.2222
$ws.Range("M31").Value = -2588.3276
$ws.Range("N31").Value = -3848.2222
$ws.Range("H34").Value = 2972.1184
$ws.Range("I34").Value = 2883.3276
$ws.Range("J34").Value = 3258.2222
$ws.Range("K34").Value = 2883.3276
$ws.Range("L34").Value = 3258.2222
$ws.Range("M34").Value = -2681.3276
$ws.Range("N34").Value = -3662.2222
$ws.Range("H58").Value = 1759.803
$ws.Range("I58").Value = 1046.375
$ws.Range("J58").Value = 3662.2778
$ws.Range("K58").Value = 1046.375
$ws.Range("L58").Value = 3662.2778
$ws.Range("M58").Value = -843.375
$ws.Range("N58").Value = -4068.2778
$ws.Range("H99").Value = 1784.65
$ws.Range("I99").Value = 1185.2667
$ws.Range("J99").Value = 3582.8
$ws.Range("K99").Value = 1185.2667
$ws.Range("L99").Value = 3582.8
$ws.Range("M99").Value = 312.7333000000001
$ws.Range("N99").Value = -6578.8
$ws.Range("H126").Value = 1784.65
$ws.Range("I126").Value = 1185.2667
$ws.Range("J126").Value = 3582.8
$ws.Range("K126").Value = 3555.800099999999
$ws.Range("L126").Value = 10748.4
$ws.Range("M126").Value = -1085.800099999999
$ws.Range("N126").Value = -15688.4
$ws.Range("H132").Value = 4079.0576
$ws.Range("I132").Value = 4460.1514
$ws.Range("J132").Value = 3417.158
$ws.Range("K132").Value = 13380.4542
$ws.Range("L132").Value = 10251.474
$ws.Range("M132").Value = -10850.4542
$ws.Range("N132").Value = -15311.474
$ws.Range("H136").Value = 1759.803
$ws.Range("I136").Value = 1046.375
$ws.Range("J136").Value = 3662.2778
$ws.Range("K136").Value = 3139.125
$ws.Range("L136").Value = 10986.8334
$ws.Range("M136").Value = -589.125
$ws.Range("N136").Value = -16086.8334
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7576651
$ws.Range("I5").Value = 509.54544
$ws.Range("J5").Value = 15152792
$ws.Range("K5").Value = 1528.63632
$ws.Range("L5").Value = 45458376
$ws.Range("M5").Value = -1416.63632
$ws.Range("N5").Value = -45458600
$ws.Range("H36").Value = 1643.75
$ws.Range("I36").Value = 572.5
$ws.Range("K36").Value = 1717.5
$ws.Range("M36").Value = -1548.5
$ws.Range("H80").Value = 2630
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 2537.5
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 7612.5
$ws.Range("M80").Value = -8064
$ws.Range("N80").Value = -9484.5
$ws.Range("H83").Value = 2630
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2537.5
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 22837.5
$ws.Range("M83").Value = -22320
$ws.Range("N83").Value = -32197.5
$ws.Range("H122").Value = 913.88464
$ws.Range("I122").Value = 646.9
$ws.Range("J122").Value = 1080.75
$ws.Range("K122").Value = 5822.099999999999
$ws.Range("L122").Value = 9726.75
$ws.Range("M122").Value = -3372.099999999999
$ws.Range("N122").Value = -14626.75
$ws.Range("H131").Value = 29246.676
$ws.Range("J131").Value = 48779.85
$ws.Range("L131").Value = 146339.55
$ws.Range("N131").Value = -156419.55
$ws.Range("H135").Value = 7576651
$ws.Range("I135").Value = 509.54544
$ws.Range("J135").Value = 15152792
$ws.Range("K135").Value = 4585.90896
$ws.Range("L135").Value = 136375128
$ws.Range("M135").Value = -2050.90896
$ws.Range("N135").Value = -136380198
$ws.Range("H139").Value = 1533320.6
$ws.Range("I139").Value = 2430670.5
$ws.Range("J139").Value = 2547.4707
$ws.Range("K139").Value = 7292011.5
$ws.Range("L139").Value = 7642.4121
$ws.Range("M139").Value = -7286871.5
$ws.Range("N139").Value = -17922.4121
$ws.Range("H140").Value = 2674.2334
$ws.Range("I140").Value = 2295.8635
$ws.Range("J140").Value = 3714.75
$ws.Range("K140").Value = 6887.5905
$ws.Range("L140").Value = 11144.25
$ws.Range("M140").Value = -1707.5905
$ws.Range("N140").Value = -21504.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3471.7827
$ws.Range("I102").Value = 3087.1765
$ws.Range("J102").Value = 4561.5
$ws.Range("K102").Value = 3087.1765
$ws.Range("L102").Value = 4561.5
$ws.Range("M102").Value = -1465.1765
$ws.Range("N102").Value = -7805.5
$ws.Range("H113").Value = 2046.3549
$ws.Range("I113").Value = 1966.5714
$ws.Range("J113").Value = 2112.0588
$ws.Range("K113").Value = 1966.5714
$ws.Range("L113").Value = 2112.0588
$ws.Range("M113").Value = 203.4286
$ws.Range("N113").Value = -6452.0588
$ws.Range("H132").Value = 3205.2712
$ws.Range("I132").Value = 1501.3778
$ws.Range("J132").Value = 8682.071
$ws.Range("K132").Value = 4504.1334
$ws.Range("L132").Value = 26046.213
$ws.Range("M132").Value = -1974.1334
$ws.Range("N132").Value = -31106.213
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 35000
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 50000
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 50000
$ws.Range("M74").Value = -19002
$ws.Range("N74").Value = -51996
$ws.Range("H76").Value = 16666.334
$ws.Range("J76").Value = 16666.334
$ws.Range("L76").Value = 16666.334
$ws.Range("N76").Value = -17342.334
$ws.Range("H77").Value = 35000
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 50000
$ws.Range("K77").Value = 60000
$ws.Range("L77").Value = 150000
$ws.Range("M77").Value = -55008
$ws.Range("N77").Value = -159984
$ws.Range("H79").Value = 16666.334
$ws.Range("J79").Value = 16666.334
$ws.Range("L79").Value = 16666.334
$ws.Range("N79").Value = -19006.334
$ws.Range("H132").Value = 4635.087
$ws.Range("I132").Value = 4598.2
$ws.Range("J132").Value = 4752.4546
$ws.Range("K132").Value = 13794.6
$ws.Range("L132").Value = 14257.3638
$ws.Range("M132").Value = -11264.6
$ws.Range("N132").Value = -19317.3638
$ws.Range("H136").Value = 3055.5066
$ws.Range("I136").Value = 1614.9016
$ws.Range("K136").Value = 4844.7048
$ws.Range("M136").Value = -2294.7048
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 7037.1763
$ws.Range("I113").Value = 11676.223
$ws.Range("J113").Value = 1818.25
$ws.Range("K113").Value = 35028.669
$ws.Range("L113").Value = 5454.75
$ws.Range("M113").Value = -32858.669
$ws.Range("N113").Value = -9794.75
$ws.Range("H132").Value = 1838.975
$ws.Range("I132").Value = 914.6818
$ws.Range("J132").Value = 2968.6667
$ws.Range("K132").Value = 2744.0454
$ws.Range("L132").Value = 8906.000100000001
$ws.Range("M132").Value = -214.0454
$ws.Range("N132").Value = -13966.0001
